# issue #5: stock data output to json file
#
# This edit targets the "股票" (stock) worksheet:
#   1. Insert a new "property_category" column (with value "stock" for every
#      data row) right after the "total" column and before the "date" column.
#   2. Clean up a couple of pre-existing data-quality glitches in that same
#      sheet:
#        - "勤美 ." -> "勤美." (stray space before the trailing period)
#        - full/half-width comma thousand separators baked into numbers
#          stored as text (e.g. "257，006") -> plain digits ("257006")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# --- 1. Insert the new property_category column -----------------------
# Before: ... F=currency G=total H=date      I=legislator_name J=legislator_id
# After:  ... F=currency G=total H=property_category I=date J=legislator_name K=legislator_id
$ws.Columns.Item(8).Insert()

$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H27").Value = "stock"

# --- 2. Fix "勤美 ." -> "勤美." ----------------------------------------
$ws.Range("B6").Value = "勤美."

# --- 3. Strip stray full/half-width commas out of numeric text values --
$ws.Range("D19").Value = "'257006"
$ws.Range("D20").Value = "'28700"
$ws.Range("D21").Value = "'22000"
$ws.Range("D23").Value = "'21000"
$ws.Range("D25").Value = "'10000"
$ws.Range("G16").Value = "'1000000"
$ws.Range("G18").Value = "'1000000"
$ws.Range("G20").Value = "'287000"
$ws.Range("G21").Value = "'220000"
